# Rename the header cells in row 1: replace underscore-joined header
# text with space-separated text (Student_Name -> Student Name, etc.),
# keeping the same column order (A=Student Name, B=Study Hours,
# C=Sleep Hours, D=Social Media Hours).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Student Name"
$ws.Range("B1").Value = "Study Hours"
$ws.Range("C1").Value = "Sleep Hours"
$ws.Range("D1").Value = "Social Media Hours"

# Move the active selection to F1, matching the saved workbook view state.
$ws.Range("F1").Select() | Out-Null
